$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 1064
$ws.Range("K3").Value = 8184
$ws.Range("L3").Value = 1070
$ws.Range("J4").Value = 1855
$ws.Range("K4").Value = 1743
$ws.Range("L4").Value = 296
$ws.Range("L5").Value = 73
$ws.Range("K6").Value = 9120
$ws.Range("L6").Value = 1096
$ws.Range("J7").Value = 29326
$ws.Range("K7").Value = 27533
$ws.Range("L7").Value = 3599

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 8

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 66
$ws.Range("K4").Value = 104
$ws.Range("L4").Value = 17
$ws.Range("L5").Value = 7
$ws.Range("L6").Value = 67
$ws.Range("K7").Value = 1803
$ws.Range("L7").Value = 224

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 40
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 21
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 115
$ws.Range("K8").Value = 1803
$ws.Range("L8").Value = 224
$ws.Range("K14").Value = 128
$ws.Range("L19").Value = 111
$ws.Range("L20").Value = 93
$ws.Range("L23").Value = 38
$ws.Range("K29").Value = 1520
$ws.Range("L29").Value = 167
$ws.Range("L36").Value = 62
$ws.Range("L37").Value = 123
$ws.Range("L42").Value = 114
$ws.Range("L43").Value = 26
$ws.Range("L44").Value = 24
$ws.Range("L51").Value = 48
$ws.Range("L53").Value = 45
$ws.Range("L54").Value = 77
$ws.Range("J63").Value = 204
$ws.Range("L65").Value = 73
$ws.Range("L67").Value = 132
$ws.Range("L69").Value = 8
$ws.Range("L79").Value = 103
$ws.Range("L80").Value = 15
$ws.Range("L84").Value = 38
$ws.Range("L87").Value = 13
$ws.Range("L88").Value = 54
$ws.Range("L91").Value = 48
$ws.Range("L93").Value = 20
$ws.Range("L95").Value = 54
$ws.Range("L97").Value = 43
$ws.Range("J101").Value = 29326
$ws.Range("K101").Value = 27533
$ws.Range("L101").Value = 3599

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 44
$ws.Range("L3").Value = 33
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 56
$ws.Range("K6").Value = 448
$ws.Range("K7").Value = 1520
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 32
$ws.Range("L3").Value = 36
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L3").Value = 13
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 30
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 29
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 27
$ws.Range("L5").Value = 5
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L3").Value = 12
$ws.Range("L4").Value = 9

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 13
